$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear A9 (used to hold "WATER" duplicate label above the second table)
$ws.Range("A9").Value = ""

# Shift the element-name column of the second table down by one row
# and append a new row (15) for the sixth Hindu god (Agni)
$ws.Range("A10").Value = "WATER"
$ws.Range("A11").Value = "WAR"
$ws.Range("A12").Value = "LOVE"
$ws.Range("A13").Value = "KING"
$ws.Range("A14").Value = "DEATH"
$ws.Range("A15").Value = "MESSENGER"

# Expand the Hindu god names into full ability descriptions
$ws.Range("B10").Value = "Veruna: can use a spell card free of cost"
$ws.Range("B11").Value = "Karttikeyn: half of your creatures can attack again this turn"
$ws.Range("B12").Value = "Rati: Can heal one of your creatures"
$ws.Range("B13").Value = "Indra: able to look at a persons hand"
$ws.Range("B14").Value = "Yama: either take the top card from your discard pile or your opponents"
$ws.Range("B15").Value = "Agni: Draw an extra card or two"

# Keep the NEUMONT ability description aligned on the new row 15
$ws.Range("C15").Value = "Take on too much: Double mana"

# Update the saved selection to match the new last-used cell
$ws.Range("B15").Select()
